$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fill TALLER 1 / TALLER 2 AREA / TALLER 3 AREA GUARDANDO with "5.0"
# (copy the existing "5.0" text cell so the written cell stays plain text,
# same as the other rows, without introducing a new number format/style)
$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("D4").Copy()
$ws.Range("E3").PasteSpecial(-4163)
$ws.Range("D4").Copy()
$ws.Range("F3").PasteSpecial(-4163)

# Bump OROS for this row from 10 to 30
$ws.Range("G3").Value = 30

# Move the active selection to A5
$ws.Range("A5").Select()
